$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values re-pulled from source data.
$updates = @{
    2  = -5
    3  = -1
    4  = -9
    6  = -7
    7  = 0
    8  = 3
    9  = -7
    11 = 2
    12 = 5
    13 = -3
    14 = -8
    15 = 3
    16 = 2
    17 = -1
    18 = -1
    19 = -2
    20 = -2
    21 = 1
    22 = 11
    23 = -4
    24 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
